$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transit_simple")

# Mapping from old "Mode Simple" value -> new "Mode Simple" value, per Primary Mode row.
# Rows 2-34 correspond to sheet rows; column A = Primary Mode, column B = Mode Simple.
# Ordered so that new shared-string values are first introduced in the same order
# as in the target workbook: Drive, Other, Transit, Bike, Walk.
$newValues = [ordered]@{
    4  = "Drive"
    24 = "Other"
    17 = "Transit"
    3  = "Bike"
    2  = "Walk"
    5  = "Drive"
    6  = "Drive"
    7  = "Drive"
    8  = "Drive"
    9  = "Drive"
    10 = "Drive"
    11 = "Drive"
    12 = "Drive"
    13 = "Drive"
    14 = "Drive"
    15 = "Drive"
    16 = "Drive"
    18 = "Drive"
    19 = "Transit"
    20 = "Transit"
    21 = "Transit"
    22 = "Transit"
    23 = "Transit"
    25 = "Transit"
    26 = "Drive"
    27 = "Drive"
    28 = "Other"
    29 = "Other"
    30 = "Transit"
    31 = "Transit"
    32 = "Other"
    33 = "Transit"
    34 = "Other"
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

$ws.Range("B2").Select()
